# Update metrics in row 3 (the data row for the second file) as per the
# "Correcting Relevance Markers Walker (2018) - Wolters (2018)" correction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 0.7892540427751695
$ws.Range("I3").Value = 0.02912519506405969
$ws.Range("K3").Value = 276.0990990990991

$ws.Range("Q3").Value = 4
$ws.Range("R3").Value = 14
$ws.Range("S3").Value = 36
$ws.Range("T3").Value = 238
$ws.Range("U3").Value = 667
$ws.Range("V3").Value = 7553
$ws.Range("W3").Value = 7543
$ws.Range("X3").Value = 7521
$ws.Range("Y3").Value = 7319
$ws.Range("Z3").Value = 6890

$ws.Range("AF3").Value = 0.999471
$ws.Range("AG3").Value = 0.998147
$ws.Range("AH3").Value = 0.995236
$ws.Range("AI3").Value = 0.968506
$ws.Range("AJ3").Value = 0.911737
